$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    # Force the cell to hold a genuine text value (not a number) without
    # altering its number-format style: write a text-literal formula,
    # then paste-special just the value over itself. This collapses the
    # formula away while leaving the cell's stored type as text.
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# H2 / H3: numbers -> text
Set-TextValue $ws.Cells.Item(2, 8) "15"
Set-TextValue $ws.Cells.Item(3, 8) "16"

# Row 4: replace contents entirely
$ws.Range("B4").Value = "Mel"
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = "Cachorro"
$ws.Range("E4").Value = "Cão-pelado-peruano"
$ws.Range("F4").Value = "Caramelo"
$ws.Range("G4").Value = "Médio"
$ws.Range("L4").Value = 0

Set-TextValue $ws.Cells.Item(4, 1) "34"
Set-TextValue $ws.Cells.Item(4, 8) "1647894877"
Set-TextValue $ws.Cells.Item(4, 9) "2024-06-20"
Set-TextValue $ws.Cells.Item(4, 10) "2024-06-20"
Set-TextValue $ws.Cells.Item(4, 11) "Históricos/34.txt"

# Rows 5-8 are no longer part of the data: delete them entirely
$ws.Range("A5:L8").EntireRow.Delete()
